# Correction in SA algorithm and 746 logs
# Update the "Fitness" column (C) values for run_27 log data according to the
# corrected generation ranges:
#   Generation 0-12   -> 7769
#   Generation 13-37  -> 7312
#   Generation 38-51  -> 7295
#   Generation 52-250 -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 252
}

for ($row = 2; $row -le $lastRow; $row++) {
    $gen = $ws.Cells.Item($row, 2).Value()

    if ($gen -le 12) {
        $newValue = 7769
    } elseif ($gen -le 37) {
        $newValue = 7312
    } elseif ($gen -le 51) {
        $newValue = 7295
    } else {
        $newValue = 7293
    }

    $ws.Cells.Item($row, 3).Value = $newValue
}
